# Commit: "Raw and clean Data from SSA for June 13th"
# Adds the June 13, 2020 (serial 43995) observations to each tracking
# sheet of the workbook and recomputes the dependent totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. out_vars — brand-new row 14 (headers + 9 metric columns).
#    Insert the row first so it inherits row 13's per-column number
#    formats (date / integer / percentage styles), then fill in values.
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Rows.Item(14).Insert(-4121)  # xlShiftDown

$wsOut.Range("A14").Value = 43995
$wsOut.Range("B14").Value = 142690
$wsOut.Range("C14").Value = 202139
$wsOut.Range("D14").Value = 56926
$wsOut.Range("E14").Value = 16872
$wsOut.Range("F14").Value = 32.661714205620576
$wsOut.Range("G14").Value = 46605
$wsOut.Range("H14").Value = 4248
$wsOut.Range("I14").Value = 4426
$wsOut.Range("J14").Value = 401755

# ---------------------------------------------------------------------
# 2. dates_dx — row 14 cells already exist (blank); just populate them.
# ---------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Range("A14").Value = 43995
$wsDx.Range("B14").Value = 0
$wsDx.Range("C14").Value = 1
$wsDx.Range("D14").Value = 1
$wsDx.Range("E14").Value = 1
$wsDx.Range("F14").Value = 0
$wsDx.Range("G14").Value = 0
$wsDx.Range("H14").Value = 0
$wsDx.Range("I14").Value = 4

# ---------------------------------------------------------------------
# 3. dates_sx — row 14 is brand new. A14 needs the same date style as
#    A13 (s=48); copy that single cell's format, then write values.
# ---------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Range("A13").Copy()
$wsSx.Range("A14").PasteSpecial(-4122)  # xlPasteFormats

$wsSx.Range("A14").Value = 43995
$wsSx.Range("B14").Value = 0
$wsSx.Range("C14").Value = 1
$wsSx.Range("D14").Value = 0
$wsSx.Range("E14").Value = 1
$wsSx.Range("F14").Value = 1
$wsSx.Range("G14").Value = 1
$wsSx.Range("H14").Value = 0
$wsSx.Range("I14").Value = 1
$wsSx.Range("J14").Value = 1
$wsSx.Range("K14").Value = 0
$wsSx.Range("L14").Value = 0

# ---------------------------------------------------------------------
# 4. dates_deaths — A14 already exists (blank, dated style); B14:H14
#    are new plain cells.
# ---------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Range("A14").Value = 43995
$wsDeaths.Range("B14").Value = 0
$wsDeaths.Range("C14").Value = 0
$wsDeaths.Range("D14").Value = 2
$wsDeaths.Range("E14").Value = 1
$wsDeaths.Range("F14").Value = 1
$wsDeaths.Range("G14").Value = 1
$wsDeaths.Range("H14").Value = 2

# ---------------------------------------------------------------------
# 5. control_obs — new column N (day 43995) for the existing metric
#    rows, plus the day total in N20 (sum of N2:N18, mirroring the
#    existing shared SUM formula used across columns C:M).
# ---------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item("control_obs")
$wsControl.Range("N1").Value = 43995
$wsControl.Range("N2").Value = 3433
$wsControl.Range("N3").Value = 3248
$wsControl.Range("N4").Value = 3248
$wsControl.Range("N5").Value = 3248
$wsControl.Range("N6").Value = 3248
$wsControl.Range("N7").Value = 2461
$wsControl.Range("N8").Value = 5070
$wsControl.Range("N10").Value = 153
$wsControl.Range("N11").Value = 153
$wsControl.Range("N12").Value = 153
$wsControl.Range("N13").Value = 153
$wsControl.Range("N14").Value = 153
$wsControl.Range("N15").Value = 100
$wsControl.Range("N16").Value = 165
$wsControl.Range("N18").Value = 800
$wsControl.Range("N20").Formula = "=SUM(N2:N18)"

# ---------------------------------------------------------------------
# 6. View state — reproduce the author's final selections on each
#    sheet, ending on control_obs (the tab left active/selected).
# ---------------------------------------------------------------------
$wsOut.Activate()
$wsOut.Range("A14").Select()

$wsDx.Activate()
$wsDx.Range("J14").Select()

$wsSx.Activate()
$wsSx.Range("K22").Select()

$wsDeaths.Activate()
$wsDeaths.Range("E20").Select()

$wsControl.Activate()
$wsControl.Range("O10").Select()
